$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 9900
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()

$ws.Range("H23").Value = 9900
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()

$ws.Range("H29").Value = 1473.1111
$ws.Range("J29").Value = 2600
$ws.Range("L29").Value = 7800
$ws.Range("N29").Value = -8362

$ws.Range("H33").Value = 766.4400000000001
$ws.Range("I33").Value = 756.75
$ws.Range("K33").Value = 756.75
$ws.Range("M33").Value = -527.75

$ws.Range("H38").Value = 2606.3845
$ws.Range("I38").Value = 93.25
$ws.Range("J38").Value = 3723.3333
$ws.Range("K38").Value = 279.75
$ws.Range("L38").Value = 11169.9999
$ws.Range("M38").Value = 92.25
$ws.Range("N38").Value = -11913.9999

$ws.Range("H40").Value = 1532.6666
$ws.Range("I40").Value = 1534.5883
$ws.Range("J40").Value = 1528
$ws.Range("K40").Value = 1534.5883
$ws.Range("L40").Value = 1528
$ws.Range("M40").Value = -1359.5883
$ws.Range("N40").Value = -1878

$ws.Range("H58").Value = 1000
$ws.Range("I58").Value = 1000
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 3000
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -2850
$ws.Range("N58").ClearContents()

$ws.Range("H87").Value = 25053.572
$ws.Range("J87").Value = 25053.572
$ws.Range("L87").Value = 25053.572
$ws.Range("N87").Value = -27549.572

$ws.Range("H90").Value = 25053.572
$ws.Range("J90").Value = 25053.572
$ws.Range("L90").Value = 75160.716
$ws.Range("N90").Value = -87640.716

$ws.Range("H112").Value = 1389.3158
$ws.Range("I112").Value = 1120
$ws.Range("J112").Value = 1485.5
$ws.Range("K112").Value = 3360
$ws.Range("L112").Value = 4456.5
$ws.Range("M112").Value = -2252
$ws.Range("N112").Value = -6672.5

$ws.Range("H127").Value = 999.6667
$ws.Range("I127").Value = 999.6667
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 2999.0001
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = 1960.9999
$ws.Range("N127").ClearContents()

$ws.Range("H129").Value = 906.8333
$ws.Range("I129").Value = 306.125
$ws.Range("J129").Value = 1125.2727
$ws.Range("K129").Value = 918.375
$ws.Range("L129").Value = 3375.8181
$ws.Range("M129").Value = 4081.625
$ws.Range("N129").Value = -13375.8181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16356.21
$ws.Range("I32").Value = 5360.61
$ws.Range("K32").Value = 5360.61
$ws.Range("M32").Value = -5073.61

$ws.Range("H63").Value = 2659.5
$ws.Range("I63").Value = 1697.5834
$ws.Range("J63").Value = 4583.3335
$ws.Range("K63").Value = 1697.5834
$ws.Range("L63").Value = 4583.3335
$ws.Range("M63").Value = -1011.5834
$ws.Range("N63").Value = -5955.3335

$ws.Range("H66").Value = 2659.5
$ws.Range("I66").Value = 1697.5834
$ws.Range("J66").Value = 4583.3335
$ws.Range("K66").Value = 8487.916999999999
$ws.Range("L66").Value = 22916.6675
$ws.Range("M66").Value = -5055.916999999999
$ws.Range("N66").Value = -29780.6675

$ws.Range("H102").Value = 1405.1578
$ws.Range("I102").Value = 1346.5333
$ws.Range("K102").Value = 1346.5333
$ws.Range("M102").Value = 275.4666999999999

$ws.Range("H105").Value = 50000
$ws.Range("J105").Value = 50000
$ws.Range("L105").Value = 50000
$ws.Range("N105").Value = -56988

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 64
$ws.Range("I22").Value = 67.5
$ws.Range("J22").Value = 53.5
$ws.Range("K22").Value = 67.5
$ws.Range("L22").Value = 53.5
$ws.Range("M22").Value = 105.5
$ws.Range("N22").Value = -399.5

$ws.Range("H99").Value = 1512.8182
$ws.Range("I99").Value = 1490
$ws.Range("K99").Value = 1490
$ws.Range("M99").Value = 8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 222
$ws.Range("I22").Value = 221.25
$ws.Range("J22").Value = 225
$ws.Range("K22").Value = 221.25
$ws.Range("L22").Value = 225
$ws.Range("M22").Value = 128.75
$ws.Range("N22").Value = -925

$ws.Range("H62").Value = 43481004
$ws.Range("J62").Value = 76926080
$ws.Range("L62").Value = 76926080
$ws.Range("N62").Value = -76927328

$ws.Range("H65").Value = 43481004
$ws.Range("J65").Value = 76926080
$ws.Range("L65").Value = 384630400
$ws.Range("N65").Value = -384636640

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 134.14285
$ws.Range("I2").Value = 67.8
$ws.Range("J2").Value = 300
$ws.Range("K2").Value = 406.8
$ws.Range("L2").Value = 1800
$ws.Range("M2").Value = -293.8
$ws.Range("N2").Value = -2026

$ws.Range("H38").Value = 265.55554
$ws.Range("J38").Value = 295.25
$ws.Range("L38").Value = 885.75
$ws.Range("N38").Value = -1579.75

$ws.Range("H113").Value = 2962.2
$ws.Range("I113").Value = 1000
$ws.Range("K113").Value = 3000
$ws.Range("M113").Value = -830

$ws.Range("H121").Value = 2000.1111
$ws.Range("I121").Value = 271.5
$ws.Range("J121").Value = 2864.4167
$ws.Range("K121").Value = 814.5
$ws.Range("L121").Value = 8593.250100000001
$ws.Range("M121").Value = 495.5
$ws.Range("N121").Value = -11213.2501

$ws.Range("H131").Value = 939975.4
$ws.Range("J131").Value = 1092.6066
$ws.Range("L131").Value = 3277.8198
$ws.Range("N131").Value = -13357.8198

$ws.Range("H140").Value = 972.7222
$ws.Range("I140").Value = 735.82355
$ws.Range("K140").Value = 2207.47065
$ws.Range("M140").Value = 2972.52935

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 6128.2
$ws.Range("J23").Value = 6128.2
$ws.Range("L23").Value = 6128.2
$ws.Range("N23").Value = -6574.2

$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()

$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()

$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()

$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("N29").ClearContents()

$ws.Range("H80").Value = 2557.1428
$ws.Range("I80").Value = 2460
$ws.Range("J80").Value = 2611.111
$ws.Range("K80").Value = 2460
$ws.Range("L80").Value = 2611.111
$ws.Range("M80").Value = -1462
$ws.Range("N80").Value = -4607.111

$ws.Range("H83").Value = 2557.1428
$ws.Range("I83").Value = 2460
$ws.Range("J83").Value = 2611.111
$ws.Range("K83").Value = 12300
$ws.Range("L83").Value = 13055.555
$ws.Range("M83").Value = -7308
$ws.Range("N83").Value = -23039.555

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 773.0769
$ws.Range("I16").Value = 866
$ws.Range("J16").Value = 463.33334
$ws.Range("K16").Value = 866
$ws.Range("L16").Value = 463.33334
$ws.Range("M16").Value = -696
$ws.Range("N16").Value = -803.33334

$ws.Range("H22").Value = 1449807.9
$ws.Range("I22").Value = 2778002.8
$ws.Range("J22").Value = 868.0909
$ws.Range("K22").Value = 2778002.8
$ws.Range("L22").Value = 868.0909
$ws.Range("M22").Value = -2777707.8
$ws.Range("N22").Value = -1458.0909

$ws.Range("H27").Value = 1449807.9
$ws.Range("I27").Value = 2778002.8
$ws.Range("J27").Value = 868.0909
$ws.Range("K27").Value = 2778002.8
$ws.Range("L27").Value = 868.0909
$ws.Range("M27").Value = -2777895.8
$ws.Range("N27").Value = -1082.0909

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

$ws.Range("H81").Value = 52634828
$ws.Range("I81").Value = 142859800
$ws.Range("J81").Value = 3584.6667
$ws.Range("K81").Value = 285719600
$ws.Range("L81").Value = 7169.3334
$ws.Range("M81").Value = -285718539
$ws.Range("N81").Value = -9291.3334

$ws.Range("H84").Value = 52634828
$ws.Range("I84").Value = 142859800
$ws.Range("J84").Value = 3584.6667
$ws.Range("K84").Value = 1428598000
$ws.Range("L84").Value = 35846.667
$ws.Range("M84").Value = -1428592696
$ws.Range("N84").Value = -46454.667
